$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Optimizer Disabled"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Optimizer Disabled")

$ws1.Range("D4").Value = 124478
$ws1.Range("E4").Value = 1162687
$ws1.Range("F4").Value = 11544932

$ws1.Range("D5").Value = 66454
$ws1.Range("E5").Value = 92009

$ws1.Range("D6").Value = 66581
$ws1.Range("E6").Value = 87281

$ws1.Range("E7").Value = 437126
$ws1.Range("F7").Value = 553001

$ws1.Range("E8").Value = 457318
$ws1.Range("F8").Value = 481180

$ws1.Range("F9").Value = 3036502

$ws1.Range("F10").Value = 2666102

$ws1.Range("D12").Select()

# ---------------------------------------------------------------------------
# Sheet "Runs 200 - Optimizer Enabled"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Runs 200 - Optimizer Enabled")

$ws2.Range("D4").Value = 123493
$ws2.Range("E4").Value = 1154448
$ws2.Range("F4").Value = 11464153

$ws2.Range("D5").Value = 65436
$ws2.Range("E5").Value = 90830

$ws2.Range("D6").Value = 65417
$ws2.Range("E6").Value = 86117

$ws2.Range("E7").Value = 421970
$ws2.Range("F7").Value = 537040

$ws2.Range("E8").Value = 440114
$ws2.Range("F8").Value = 464379

$ws2.Range("F9").Value = 2837301

$ws2.Range("F10").Value = 2440701

$ws2.Range("H14").Select()

# ---------------------------------------------------------------------------
# Sheet "Runs 1000 - Optimizer Enabled"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Runs 1000 - Optimizer Enabled")

$ws3.Range("D4").Value = 123481
$ws3.Range("E4").Value = 1154328
$ws3.Range("F4").Value = 11462976

$ws3.Range("D5").Value = 65380
$ws3.Range("E5").Value = 90774

$ws3.Range("D6").Value = 65428
$ws3.Range("E6").Value = 86128

$ws3.Range("E7").Value = 421806
$ws3.Range("F7").Value = 536876

$ws3.Range("E8").Value = 439950
$ws3.Range("F8").Value = 464215

$ws3.Range("F9").Value = 2836657

$ws3.Range("F10").Value = 2440101

$ws3.Range("G7").Select()

# ---------------------------------------------------------------------------
# Restore the originally active sheet/tab ("Runs 200 - Optimizer Enabled")
# and its selection, since selecting ranges on other sheets above would
# otherwise shift the workbook's active tab.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("H14").Select()
